$wb = $excel.ActiveWorkbook

# Original sheets: description, license, macro (in that tab order).
$wsDescription = $wb.Worksheets.Item("description")
$wsMacro       = $wb.Worksheets.Item("macro")

# Build the new "functions" sheet by duplicating the "macro" sheet's layout
# (same table structure: name / help / body / location) and placing the
# duplicate right after "description".
$wsMacro.Copy($null, $wsDescription)

# Worksheet indices shift once the copy is inserted, so re-resolve sheets by
# name (stable) rather than trusting previously captured index-bound handles.
$wsLicense = $wb.Worksheets.Item("license")
$wsLicense.Delete()

# The duplicate got an auto-generated name ("macro (2)") - rename it.
$wsFunctions = $wb.Worksheets.Item("macro (2)")
$wsFunctions.Name = "functions"

# --- Populate the functions sheet -----------------------------------------
$wsFunctions.Range("A2").Value = "hellofun"
$wsFunctions.Range("B2").Value = "This is a hello function example."
$wsFunctions.Range("C2").Value = "function hellofun(str `$) `$ 128;`n  return(catx(`" `", `"Hello`", str, `"!!`"));`nendfunc;"

# Row 2 holds the single data row, sized a bit shorter than the source macro rows.
$wsFunctions.Rows.Item(2).RowHeight = 45

# Row 3 stays as an empty styled placeholder row (no custom height).
$wsFunctions.Rows.Item(3).ClearContents()
$wsFunctions.Rows.Item(3).EntireRow.AutoFit()

# Row 4 is removed entirely (no formatting left behind).
$wsFunctions.Rows.Item(4).ClearFormats()
$wsFunctions.Rows.Item(4).ClearContents()

# --- Update the macro sheet -------------------------------------------------
$wsMacro = $wb.Worksheets.Item("macro")
$wsMacro.Range("B3").Value = "This is macro to say hello to something.`n(e.g. ``%myhello2(obj=Taro)``)"
$wsMacro.Range("B27").Select()

# Selection / active cell on the functions sheet - selected last so it ends
# up the active tab in the saved workbook.
$wsFunctions.Range("C2").Select()
$wsFunctions.Activate()
